$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 13-14, shifting existing rows 13-67 down to 15-69
$ws.Range("13:14").Insert()

# Populate new row 13: Santina / Primera, week of 2021-12-14 (serial 44544)
$ws.Cells.Item(13, 1).Value = 7
$ws.Cells.Item(13, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(13, 3).Value = "Ñuble"
$ws.Cells.Item(13, 4).Value = 44544
$ws.Cells.Item(13, 5).Value = 16
$ws.Cells.Item(13, 6).Value = "Fruta"
$ws.Cells.Item(13, 7).Value = 100103
$ws.Cells.Item(13, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(13, 9).Value = 100103001
$ws.Cells.Item(13, 10).Value = "Cereza"
$ws.Cells.Item(13, 11).Value = "Santina"
$ws.Cells.Item(13, 12).Value = "Primera"
$ws.Cells.Item(13, 13).Value = 120
$ws.Cells.Item(13, 14).Value = 8500
$ws.Cells.Item(13, 15).Value = 9000
$ws.Cells.Item(13, 16).Value = 8750
$ws.Cells.Item(13, 17).Value = "$/caja 10 kilos"
$ws.Cells.Item(13, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(13, 19).Value = 875
$ws.Cells.Item(13, 20).Value = 10

$ws.Cells.Item(14, 1).Value = 7
$ws.Cells.Item(14, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(14, 3).Value = "Ñuble"
$ws.Cells.Item(14, 4).Value = 44544
$ws.Cells.Item(14, 5).Value = 16
$ws.Cells.Item(14, 6).Value = "Fruta"
$ws.Cells.Item(14, 7).Value = 100103
$ws.Cells.Item(14, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(14, 9).Value = 100103001
$ws.Cells.Item(14, 10).Value = "Cereza"
$ws.Cells.Item(14, 11).Value = "Santina"
$ws.Cells.Item(14, 12).Value = "Segunda"
$ws.Cells.Item(14, 13).Value = 60
$ws.Cells.Item(14, 14).Value = 7000
$ws.Cells.Item(14, 15).Value = 7500
$ws.Cells.Item(14, 16).Value = 7250
$ws.Cells.Item(14, 17).Value = "$/caja 10 kilos"
$ws.Cells.Item(14, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(14, 19).Value = 725
$ws.Cells.Item(14, 20).Value = 10

Write-Host "Inserted 2 rows and populated new weekly data."
